$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (engine stores widths quantized to 1/6-character pixel
# units, same as real Excel; these inputs land on the closest achievable
# stored widths to the target 15.42578125 / 14.7109375)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666668
$ws.Columns.Item(2).ColumnWidth = 13.833333333333332

# Update cell values
$ws.Range("A1").Value = -0.021144659550824819
$ws.Range("B1").Value = -0.02115841701288447

$ws.Range("A2").Value = -0.025596668085641247
$ws.Range("B2").Value = 0.025596668054237829

$ws.Range("A3").Value = -0.0051682340627409529
$ws.Range("B3").Value = 0.0051682340069879214

$ws.Range("A4").Value = 0.047966594128354183
$ws.Range("B4").Value = -0.047966594176531062

$ws.Range("A5").Value = -0.014271729862585163
$ws.Range("B5").Value = 0.014271729781320261
